$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Вид ДК"
$ws.Range("B1").Value = "Штрих-код карты"
$ws.Range("C1").Value = "Фамилия"
$ws.Range("D1").Value = "Имя"
$ws.Range("E1").Value = "Отчество"
$ws.Range("F1").Value = "Дата рождения"
$ws.Range("H1").Value = "Дом. Телефон"
$ws.Range("I1").Value = "СМС"
$ws.Range("J1").Value = "Моб. Телефон"
$ws.Range("K1").Value = "Общая сумма покупок"
$ws.Range("L1").Value = "Общая сумма бонусов"
$ws.Range("M1").Value = "Изменить"

# ---- Row 2 ----
$ws.Range("A2").Value = "пластик"
$ws.Range("B2").Value = 2701200000000

# ---- Row 3 ----
$ws.Range("A3").Value = "пластик"
$ws.Range("B3").Value = 2701200000001

# ---- Row 4 ----
$ws.Range("A4").Value = "пластик"
$ws.Range("B4").Value = 2701200000002
$ws.Range("C4").Value = "Хрушков"
$ws.Range("D4").Value = "Степан"
$ws.Range("E4").Value = "Игоревич"

# Force the birth date and phone number to stay as plain text rather than
# being auto-parsed into a date serial / number.
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "10.01.2005"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "79965677951"
$ws.Range("F4").ClearFormats()
$ws.Range("J4").ClearFormats()

# ---- Number formats ----
# Card barcode column -> integer "0" format
$ws.Range("B1:B4").NumberFormat = "0"
# Birth-date header -> "0.00" format (matches the source workbook)
$ws.Range("F1").NumberFormat = "0.00"

# ---- Column widths (closest achievable match to the source widths) ----
$ws.Range("A1").EntireColumn.ColumnWidth = 6.8333333333
$ws.Range("B1").EntireColumn.ColumnWidth = 15.6666666667
$ws.Range("C1").EntireColumn.ColumnWidth = 8.1666666667
$ws.Range("D1").EntireColumn.ColumnWidth = 3.8333333333
$ws.Range("E1").EntireColumn.ColumnWidth = 8
$ws.Range("F1").EntireColumn.ColumnWidth = 13.6666666667
$ws.Range("H1").EntireColumn.ColumnWidth = 12.5
$ws.Range("I1").EntireColumn.ColumnWidth = 4
$ws.Range("J1").EntireColumn.ColumnWidth = 12.6666666667
$ws.Range("K1").EntireColumn.ColumnWidth = 19.5
$ws.Range("L1").EntireColumn.ColumnWidth = 19.8333333333
$ws.Range("M1").EntireColumn.ColumnWidth = 8.6666666667

# ---- Selection (matches the saved workbook view) ----
$ws.Range("C7").Select() | Out-Null
